$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 350, shifting existing rows 350:425 down to 351:426
$ws.Rows.Item(350).Insert()

# Populate the newly inserted row 350 with its data
$ws.Range("A350").Value2 = 4
$ws.Range("B350").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C350").Value2 = "Los Lagos"
$ws.Range("D350").Value2 = 45244
$ws.Range("E350").Value2 = 10
$ws.Range("F350").Value2 = 100112039
$ws.Range("G350").Value2 = "Ciboulette"
$ws.Range("H350").Value2 = "Sin especificar"
$ws.Range("I350").Value2 = "Primera"
$ws.Range("J350").Value2 = 240
$ws.Range("K350").Value2 = 3500
$ws.Range("L350").Value2 = 3500
$ws.Range("M350").Value2 = 3500
$ws.Range("N350").Value2 = "$/docena de atados"
$ws.Range("O350").Value2 = "Región Metropolitana"
$ws.Range("P350").Value2 = 1167
$ws.Range("Q350").Value2 = 3
$ws.Range("R350").Value2 = "Hortaliza"
